# SCD0267 - fix test-data values, align F3's formatting with F4, and
# refresh the sheet's scroll/zoom/selection state.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$win = $excel.ActiveWindow

# F3 previously used a distinct "Calibri 11 / black" font; bring it in line
# with F4 (and the rest of the column), which uses the sheet's normal
# "Arial 9" cell style, by copying F4's formatting onto F3.
$ws.Range("F4").Copy()
[void]$ws.Range("F3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Corrected test data values.
$ws.Range("F3").Value = 32362
$ws.Range("F4").Value = 19331

# Restore default scroll position (top-left back to A1), rezoom the sheet,
# and move the active selection to F3.
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Zoom = 60
[void]$ws.Range("F3").Select()
